# update scripts wuth new tpm
#
# The underlying NATMI ligand-receptor analysis was re-run against the
# updated TPM expression matrix. Sending/target cluster membership, gene
# symbols and the sheet layout are unchanged, but every downstream
# expression-derived statistic (detection counts/rates, average/total
# expression, derived-specificity scores, and the edge-weight columns)
# was recomputed, so those cached numeric results are refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs, Bmp4-Bmpr1a / target ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.202518666666666
$ws.Range("H2").Value = 12.607556
$ws.Range("I2").Value = 0.08075097102331126
$ws.Range("J2").Value = 0.08075097102331129
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 9.839590894472444
$ws.Range("R2").Value = 88.556318050252
$ws.Range("S2").Value = 0.003208650302632945
$ws.Range("T2").Value = 0.003208650302632946

# Row 3 (target FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.202518666666666
$ws.Range("H3").Value = 12.607556
$ws.Range("I3").Value = 0.08075097102331126
$ws.Range("J3").Value = 0.08075097102331129
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("Q3").Value = 143.1074169859418
$ws.Range("R3").Value = 1287.966752873476
$ws.Range("S3").Value = 0.04666674272798416
$ws.Range("T3").Value = 0.04666674272798418

# Row 4 (target MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.202518666666666
$ws.Range("H4").Value = 12.607556
$ws.Range("I4").Value = 0.08075097102331126
$ws.Range("J4").Value = 0.08075097102331129
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("Q4").Value = 94.68250741727556
$ws.Range("R4").Value = 852.1425667554799
$ws.Range("S4").Value = 0.03087557799269416
$ws.Range("T4").Value = 0.03087557799269417

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.7496282157262072
$ws.Range("J5").Value = 0.7496282157262073
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 91.34298785794034
$ws.Range("R5").Value = 822.0868907214631
$ws.Range("S5").Value = 0.02978657433800674
$ws.Range("T5").Value = 0.02978657433800674

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.7496282157262072
$ws.Range("J6").Value = 0.7496282157262073
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("Q6").Value = 1328.496193827674
$ws.Range("S6").Value = 0.4332171692998451
$ws.Range("T6").Value = 0.4332171692998452

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.7496282157262072
$ws.Range("J7").Value = 0.7496282157262073
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("S7").Value = 0.2866244720883554
$ws.Range("T7").Value = 0.2866244720883554

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 8.827567333333333
$ws.Range("I8").Value = 0.1696208132504815
$ws.Range("J8").Value = 0.1696208132504815
$ws.Range("M8").Value = 2.341355666666667
$ws.Range("N8").Value = 7.024067000000001
$ws.Range("O8").Value = 0.03973512964576821
$ws.Range("P8").Value = 0.0397351296457682
$ws.Range("Q8").Value = 20.66847479878156
$ws.Range("R8").Value = 186.016273189034
$ws.Range("S8").Value = 0.00673990500512852
$ws.Range("T8").Value = 0.006739905005128521

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 8.827567333333333
$ws.Range("I9").Value = 0.1696208132504815
$ws.Range("J9").Value = 0.1696208132504815
$ws.Range("O9").Value = 0.5779093692199981
$ws.Range("P9").Value = 0.5779093692199981
$ws.Range("Q9").Value = 300.6031524292602
$ws.Range("R9").Value = 2705.428371863342
$ws.Range("S9").Value = 0.09802545719216885
$ws.Range("T9").Value = 0.09802545719216889

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 8.827567333333333
$ws.Range("I10").Value = 0.1696208132504815
$ws.Range("J10").Value = 0.1696208132504815
$ws.Range("O10").Value = 0.3823555011342337
$ws.Range("P10").Value = 0.3823555011342337
$ws.Range("Q10").Value = 198.8845917911844
$ws.Range("S10").Value = 0.0648554510531841
$ws.Range("T10").Value = 0.06485545105318413
